$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: experiment_description
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("experiment_description")

# Insert a new row 14 (experiment 13: "Ferguson Iceland age-simple"), pushing
# the existing 1001/1002/1003 rows down by one.
$ws1.Rows("14:14").Insert()
$ws1.Rows("14:14").RowHeight = 17

$ws1.Range("A14").Value = 13
$ws1.Range("B14").Value = "Ferguson Iceland age-simple"
$ws1.Range("C14").Value = "ferguson"
$ws1.Range("D14").Value = "age_official"
$ws1.Range("E14").Value = "iceland"
$ws1.Range("F14").Value = "iceland"

# New label_english / label_icelandic columns (G/H).
$ws1.Range("G1").Value = "label_english"
$ws1.Range("H1").Value = "label_icelandic"

$ws1.Range("G2").Value = "Base model"
$ws1.Range("H2").Value = "Grunnlíkan"

$ws1.Range("G5").Value = "Improved model"
$ws1.Range("H5").Value = "Bætt líkan"

$ws1.Range("G11").Value = "Ferguson Wuhan"
$ws1.Range("H11").Value = "Ferguson Wuhan"

$ws1.Range("G12").Value = "Ferguson Iceland"
$ws1.Range("H12").Value = "Ferguson Ísland"

$ws1.Range("G15").Value = "CA model"
$ws1.Range("H15").Value = "KM grunnlíkan"

# ---------------------------------------------------------------------------
# Sheet 2: experiment_specification
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("experiment_specification")

# Insert 3 new rows at row 38 for the new experiment 13 (Ferguson Iceland
# age-simple), pushing the existing 1001/1002/1003 blocks down by 3 rows.
$ws2.Rows("38:40").Insert()

$ws2.Range("A38").Value = 13
$ws2.Range("B38").Value = "home"
$ws2.Range("C38").Value = "length_of_stay_simple_week"
$ws2.Range("D38").Value = "none"
$ws2.Range("E38").Value = "age_simple"

$ws2.Range("A39").Value = 13
$ws2.Range("B39").Value = "inpatient_ward"
$ws2.Range("C39").Value = "length_of_stay_simple_week"
$ws2.Range("D39").Value = "none"
$ws2.Range("E39").Value = "age_simple"

$ws2.Range("A40").Value = 13
$ws2.Range("B40").Value = "intensive_care_unit"
$ws2.Range("C40").Value = "none"
$ws2.Range("D40").Value = "none"
$ws2.Range("E40").Value = "none"

# ---------------------------------------------------------------------------
# Sheet 3: run_description
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("run_description")

$ws3.Range("A10").Value = 9
$ws3.Range("B10").Value = "Compare base and best with Ferguson"

$ws3.Range("A11").Value = 10
$ws3.Range("B11").Value = "Compare Ferguson Iceland models"

# ---------------------------------------------------------------------------
# Sheet 4: run_specification
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("run_specification")

$ws4.Range("A25").Value = 9
$ws4.Range("B25").Value = 1

$ws4.Range("A26").Value = 9
$ws4.Range("B26").Value = 4

$ws4.Range("A27").Value = 9
$ws4.Range("B27").Value = 10

$ws4.Range("A28").Value = 9
$ws4.Range("B28").Value = 11

$ws4.Range("A29").Value = 10
$ws4.Range("B29").Value = 11

$ws4.Range("A30").Value = 10
$ws4.Range("B30").Value = 13

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------
[void]$ws2.Range("E39").Select()
[void]$ws3.Range("B16").Select()
[void]$ws4.Range("A31").Select()
[void]$ws1.Activate()
[void]$ws1.Range("G23").Select()
